$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Astronauta")
$ws2 = $wb.Worksheets.Item("Senador")
$ws3 = $wb.Worksheets.Item("Mago")
$ws4 = $wb.Worksheets.Item("Ninja")

# Astronauta: fill column C (T2 test results)
$ws1.Range("C2").Value = 1
$ws1.Range("C3").Value = 1
$ws1.Range("C4").Value = 0.5
$ws1.Range("C5").Value = 0.8
$ws1.Range("C6").Value = 0.8
$ws1.Range("C7").Value = 1
$ws1.Range("C8").Value = 0.8
$ws1.Range("C9").Value = 1
$ws1.Range("C10").Value = 1
$ws1.Range("C11").Value = 1
$ws1.Range("C12").Value = 1
$ws1.Range("C13").Value = 1
$ws1.Range("C14").Value = 0.5
$ws1.Range("C15").Value = 1
$ws1.Range("C16").Value = 0.5
$ws1.Range("C17").Value = 1
$ws1.Range("C18").Value = 1
$ws1.Range("C19").Value = 1
$ws1.Range("C20").Value = 1
$ws1.Range("C21").Value = 1

# Senador: fill column C (T2 test results)
$ws2.Range("C2").Value = 1
$ws2.Range("C3").Value = 1
$ws2.Range("C4").Value = 0.5
$ws2.Range("C5").Value = 0.5
$ws2.Range("C6").Value = 0.5
$ws2.Range("C7").Value = 1
$ws2.Range("C8").Value = 1
$ws2.Range("C9").Value = 1
$ws2.Range("C10").Value = 1
$ws2.Range("C11").Value = 1
$ws2.Range("C12").Value = 1
$ws2.Range("C13").Value = 0.7
$ws2.Range("C14").Value = 0.7
$ws2.Range("C15").Value = 0.7
$ws2.Range("C16").Value = 0.5
$ws2.Range("C17").Value = 1
$ws2.Range("C18").Value = 0.7
$ws2.Range("C19").Value = 1
$ws2.Range("C20").Value = 1
$ws2.Range("C21").Value = 1

# Mago: fill column D (T3 test results)
$ws3.Range("D2").Value = 1
$ws3.Range("D3").Value = 1
$ws3.Range("D4").Value = 1
$ws3.Range("D5").Value = 0.5
$ws3.Range("D6").Value = 0.8
$ws3.Range("D7").Value = 1
$ws3.Range("D8").Value = 1
$ws3.Range("D9").Value = 1
$ws3.Range("D10").Value = 1
$ws3.Range("D11").Value = 1
$ws3.Range("D12").Value = 1
$ws3.Range("D13").Value = 1
$ws3.Range("D14").Value = 0.8
$ws3.Range("D15").Value = 1
$ws3.Range("D16").Value = 0.5
$ws3.Range("D17").Value = 1
$ws3.Range("D18").Value = 0.5
$ws3.Range("D19").Value = 1
$ws3.Range("D20").Value = 1
$ws3.Range("D21").Value = 1

# Update selections per sheet and the active sheet/tab (Ninja ends up active)
$ws1.Activate()
$ws1.Range("C19").Select()

$ws2.Activate()
$ws2.Range("C20").Select()

$ws3.Activate()
$ws3.Range("D21").Select()

$ws4.Activate()
$ws4.Range("D2:D21").Select()
